$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tardigrade")

$ws.Cells.Item(1, 13).Value = "13.1.0"
$ws.Cells.Item(1, 13).Font.Bold = $true

$ws.Cells.Item(2, 13).Value = "11.1.0"

$ws.Range("M2").Select()
